# Update the regression coefficient labels to match the revised results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-0.37***"
$ws.Range("B3").Value = "-3.46***"
$ws.Range("C3").Value = "-0.81***"
